$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# max_depth (row 6) was empty -> set to 10
$ws.Range("B6").Value = 10

# min_samples_split (row 12): 5 -> 2
$ws.Range("B12").Value = 2

# n_estimators (row 15): 100 -> 200
$ws.Range("B15").Value = 200
